$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing full range (A1:A125) before writing the new, shorter amenity list
$ws.Range("A1:A125").ClearContents()

# Each value below has an extra leading apostrophe because Excel treats a single
# leading apostrophe as a text-qualifier prefix and strips it; doubling it up
# preserves the literal apostrophe that is part of each amenity label.
$values = @(
    "''Air Conditioning'",
    "''Balcony'",
    "''Basement'",
    "''Breakfast Nook'",
    "''Built-In Bookshelves'",
    "''Cable Ready'",
    "''Carpet'",
    "''Ceiling Fans'",
    "''Crown Molding'",
    "''Deck'",
    "''Den'",
    "''Dining Room'",
    "''Dishwasher'",
    "''Disposal'",
    "''Dock'",
    "''Double Pane Windows'",
    "''Double Vanities'",
    "''Eat-in Kitchen'",
    "''Family Room'",
    "''Fireplace'",
    "''Freezer'",
    "''Furnished'",
    "''Furnished Units Available'",
    "''Gated'",
    "''Granite Countertops'",
    "''Greenhouse'",
    "''Handrails'",
    "''Hardwood Floors'",
    "''Heating'",
    "''High Speed Internet Access'",
    "''Ice Maker'",
    "''In Unit Washer & Dryer'",
    "''Instant Hot Water'",
    "''Intercom'",
    "''Island Kitchen'",
    "''Kitchen'",
    "''Large Bedrooms'",
    "''Laundry Facilities'",
    "''Lawn'",
    "''Linen Closet'",
    "''Loft Layout'",
    "''Microwave'",
    "''Office'",
    "''Oven'",
    "''Pantry'",
    "''Patio'",
    "''Playground'",
    "''Porch'",
    "''Range'",
    "''Refrigerator'",
    "''Security System'",
    "''Skylight'",
    "''Smoke Free'",
    "''Sprinkler System'",
    "''Stainless Steel Appliances'",
    "''Storage Units'",
    "''Sunroom'",
    "''Tile Floors'",
    "''Trash Compactor'",
    "''Tub/Shower'",
    "''Vaulted Ceiling'",
    "''Views'",
    "''Vinyl Flooring'",
    "''Walk-In Closets'",
    "''Washer/Dryer Hookup'",
    "''Wet Bar'",
    "''Wheelchair Accessible (Rooms)'",
    "''Wi-Fi'",
    "''Window Coverings'",
    "''Yard'"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

